# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.831.40"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "2.033.91"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.56"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.36"
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.379"
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0819"
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "2.336.24"
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.50"
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.21"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.760"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.18"
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("D17").Value = "2.041.85"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("D18").Value = "37.775.16"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.83"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.90"
$ws.Range("E20").Value = "  -6.12%  "
$ws.Range("D21").Value = "0.0₃0824"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.98"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.35"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.30"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.130"
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.87"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("E30").Value = "  -4.26%  "
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.24"
$ws.Range("E32").Value = "  +8.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.41"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.50"
$ws.Range("E35").Value = "  -2.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.35"
$ws.Range("E36").Value = "  +3.49%  "
$ws.Range("E37").Value = "  -2.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.34"
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.64"
$ws.Range("E40").Value = "  +4.24%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.532.43"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0217"
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.32"
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("E44").Value = "  -2.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0914"
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("E46").Value = "  -3.32%  "
$ws.Range("E47").Value = "  -2.90%  "
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.96"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.12"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").Value = "2.225.12"
$ws.Range("E51").Value = "  -1.27%  "
